# Generate Report for Handback
# This script updates the localization-status workbook so that the two
# translation sheets (zh-cn, de-de) reflect that the handback files have
# arrived: it fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns, flips the Status text, adds the new
# hyperlinks for the target-file column, and widens a few columns so the
# new text fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$srcUrl89 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/92d704832c53f77d51874fec774ee43eb1412e26/e2e/89ade265-732c-455e-a9a8-e7c29f596ac6.md"
$srcUrlE2a = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/92d704832c53f77d51874fec774ee43eb1412e26/e2e/e2a3e336-f72f-43e5-a550-e36673991c45.md"

$md89  = "89ade265-732c-455e-a9a8-e7c29f596ac6.md"
$mdE2a = "e2a3e336-f72f-43e5-a550-e36673991c45.md"

# ---------------------------------------------------------------------
# Overview sheet: the "Status" values shown here are the same shared
# string as the per-language sheets, so just updating the text anywhere
# keeps them all in sync, but set them explicitly for clarity.
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("I2").Value = $md89
$zhcn.Range("J2").Value = "89ade265-732c-455e-a9a8-e7c29f596ac6.cfdc510da5510389a906b07e242a34dae50c84bf.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-20 05:03:28"

$zhcn.Range("I3").Value = $mdE2a
$zhcn.Range("J3").Value = "e2a3e336-f72f-43e5-a550-e36673991c45.546411cd70038f291bf2287ddd9701c283fac851.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-20 05:03:28"

# Rebuild the hyperlinks in row order (source link, target link) per row
# so relationship ids come out as rId2..rId5 in the same order Excel
# would assign them when regenerating the report.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $srcUrl89, "", "", $md89)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $srcUrl89, "", "", $md89)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $srcUrlE2a, "", "", $mdE2a)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $srcUrlE2a, "", "", $mdE2a)

$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 0xED9564
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = 0xED9564

$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(9).ColumnWidth = 39.1
$zhcn.Columns.Item(10).ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = $md89
$dede.Range("J2").Value = "89ade265-732c-455e-a9a8-e7c29f596ac6.cfdc510da5510389a906b07e242a34dae50c84bf.de-de.xlf"
$dede.Range("K2").Value = "2016-08-20 05:03:35"

$dede.Range("I3").Value = $mdE2a
$dede.Range("J3").Value = "e2a3e336-f72f-43e5-a550-e36673991c45.546411cd70038f291bf2287ddd9701c283fac851.de-de.xlf"
$dede.Range("K3").Value = "2016-08-20 05:03:35"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $srcUrl89, "", "", $md89)
$dede.Hyperlinks.Add($dede.Range("I2"), $srcUrl89, "", "", $md89)
$dede.Hyperlinks.Add($dede.Range("A3"), $srcUrlE2a, "", "", $mdE2a)
$dede.Hyperlinks.Add($dede.Range("I3"), $srcUrlE2a, "", "", $mdE2a)

$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 0xED9564
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = 0xED9564

$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(9).ColumnWidth = 39.1
$dede.Columns.Item(10).ColumnWidth = 39.1

$wb.Save()
